# Rename the "_old" / "_new" header-suffix convention to the concrete
# format-version identifiers used by this comparison (FV2304 -> FV2310),
# turn the header row into a real Excel Table (so column headers pick up
# an AutoFilter + structured references), and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the column headers: "<Name>_old" -> "<Name>_FV2304",
#    "<Name>_new" -> "<Name>_FV2310". "diff" (column K) is left as-is.
$headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304",
    "diff",
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2) Turn the populated range into an Excel Table ("Table1"), which also
#    adds the AutoFilter on the header row. The table column names are
#    picked up automatically from the (just renamed) header cells.
$dataRange = $ws.Range("A1:U85")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"

# 3) Freeze the header row so it stays visible while scrolling.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
